# Updates the cryptocurrency price/volume table in-place to match the
# latest scrape (per the "Updated cryptos list" GitHub Actions commit).
# Numeric-looking price strings (column D) are assigned with a leading
# apostrophe so Excel keeps them as text, consistent with the rest of
# the column (prices like "3.478.48" use '.' as both thousands and
# decimal separators and must never be auto-converted to a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.618.13'
$ws.Range("E2").Value = '  -4.57%  '

$ws.Range("D3").Value = '3.471.01'
$ws.Range("E3").Value = '  -5.91%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''604.69'
$ws.Range("E5").Value = '  -6.67%  '

$ws.Range("D6").Value = '''148.99'
$ws.Range("E6").Value = '  -7.68%  '

$ws.Range("D7").Value = '3.469.48'
$ws.Range("E7").Value = '  -5.89%  '

$ws.Range("E9").Value = '  -4.71%  '

$ws.Range("E10").Value = '  -5.52%  '

$ws.Range("D11").Value = '''6.89'
$ws.Range("E11").Value = '  -3.97%  '

$ws.Range("E12").Value = '  -5.20%  '

$ws.Range("E13").Value = '  -5.87%  '

$ws.Range("D14").Value = '4.052.48'
$ws.Range("E14").Value = '  -6.06%  '

$ws.Range("D15").Value = '''31.41'
$ws.Range("E15").Value = '  -4.11%  '

$ws.Range("D16").Value = '3.456.63'
$ws.Range("E16").Value = '  -6.44%  '

$ws.Range("D17").Value = '66.518.37'

$ws.Range("E18").Value = '  -0.33%  '

$ws.Range("D20").Value = '''14.99'
$ws.Range("E20").Value = '  -6.68%  '

$ws.Range("D21").Value = '''442.90'
$ws.Range("E21").Value = '  -6.25%  '

$ws.Range("D22").Value = '''9.02'
$ws.Range("E22").Value = '  -13.56%  '

$ws.Range("D23").Value = '''0.622'

$ws.Range("D24").Value = '''77.15'
$ws.Range("E24").Value = '  -3.77%  '

$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("D26").Value = '3.605.13'
$ws.Range("E26").Value = '  -6.08%  '

$ws.Range("E27").Value = '  -3.33%  '

$ws.Range("D28").Value = '''10.05'
$ws.Range("E28").Value = '  -8.50%  '

$ws.Range("D29").Value = '''8.22'
$ws.Range("E29").Value = '  -10.25%  '

$ws.Range("D30").Value = '''2.51'
$ws.Range("E30").Value = '  -5.59%  '

$ws.Range("D31").Value = '''1.56'
$ws.Range("E31").Value = '  -9.11%  '

$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  +0.12%  '

$ws.Range("E33").Value = '  -4.15%  '

$ws.Range("E34").Value = '  -4.95%  '

$ws.Range("D35").Value = '''6.12'
$ws.Range("E35").Value = '  -6.31%  '

$ws.Range("E36").Value = '  -7.96%  '

$ws.Range("D37").Value = '3.455.04'
$ws.Range("E37").Value = '  -6.38%  '

$ws.Range("D38").Value = '''7.92'
$ws.Range("E38").Value = '  -6.01%  '

$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("E40").Value = '  -0.29%  '

$ws.Range("D41").Value = '''173.23'
$ws.Range("E41").Value = '  -3.03%  '

$ws.Range("D42").Value = '''2.15'
$ws.Range("E42").Value = '  -4.14%  '

$ws.Range("D43").Value = '''0.0863'
$ws.Range("E43").Value = '  -4.70%  '

$ws.Range("E44").Value = '  -7.12%  '

$ws.Range("D45").Value = '''0.879'
$ws.Range("E45").Value = '  -5.59%  '

$ws.Range("D46").Value = '''45.30'
$ws.Range("E46").Value = '  -3.27%  '

$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '''1.23'
$ws.Range("E47").Value = '  -2.93%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '''26.31'
$ws.Range("E48").Value = '  -10.07%  '

$ws.Range("D49").Value = '''2.50'
$ws.Range("E49").Value = '  -11.39%  '

$ws.Range("E50").Value = '  -4.01%  '

$ws.Range("D51").Value = '''1.01'
$ws.Range("E51").Value = '  -4.58%  '
